# Updated cryptos list on Fri May 10 15:51:23 UTC 2024 with GitHub Actions
# Refreshes price / 1h-volume-change figures for each coin row, and swaps
# the InjectiveProtocol / USDe rows (49 and 50) to reflect the new ranking.
# Numeric-looking text in column D is written with a leading apostrophe so
# Excel keeps it as text (preserving formatting such as trailing zeros /
# thousand-dot separators) instead of silently coercing it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.844.41"
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("D3").Value = "2.931.64"
$ws.Range("E3").Value = "  -2.35%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'585.56"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").Value = "'146.66"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "2.914.28"
$ws.Range("E8").Value = "  -2.84%  "
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("D10").Value = "'6.76"
$ws.Range("E10").Value = "  +7.95%  "
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("D12").Value = "'0.449"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "'0.0000224"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("D14").Value = "'34.58"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D16").Value = "3.426.00"
$ws.Range("E16").Value = "  -2.24%  "
$ws.Range("D17").Value = "'6.85"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "60.922.63"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").Value = "2.936.51"
$ws.Range("E19").Value = "  -2.32%  "
$ws.Range("D20").Value = "'427.77"
$ws.Range("E20").Value = "  -4.57%  "
$ws.Range("D21").Value = "'13.82"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").Value = "'0.673"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").Value = "'7.23"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").Value = "'80.50"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "'10.90"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("E26").Value = "  -2.43%  "
$ws.Range("D27").Value = "'11.89"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'7.36"
$ws.Range("E29").Value = "  +2.46%  "
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").Value = "'2.19"
$ws.Range("E31").Value = "  +5.16%  "
$ws.Range("D32").Value = "'2.63"
$ws.Range("E32").Value = "  -2.94%  "
$ws.Range("D33").Value = "'26.86"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("D35").Value = "0.0₃0829"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("D37").Value = "'5.70"
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("D39").Value = "'3.00"
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("E40").Value = "  +1.68%  "
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").Value = "'8.76"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("E43").Value = "  +6.41%  "
$ws.Range("D44").Value = "'41.66"
$ws.Range("E44").Value = "  +3.22%  "
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").Value = "'372.52"
$ws.Range("E46").Value = "  -6.85%  "
$ws.Range("D47").Value = "2.657.65"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").Value = "'133.31"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'25.35"
$ws.Range("E50").Value = "  +7.49%  "
$ws.Range("E51").Value = "  -0.96%  "
